# Weekly update: a new price observation for the "Achicoria" sheet is added
# as a new row at position 24 (most recent week), pushing every existing
# row from 24 downward by one. This mirrors how the upstream dataset
# prepends the newest weekly record near the top of its date-ordered block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 24; Excel shifts rows 24..75 down to 25..76
# and grows the sheet dimension to A1:R76 automatically.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly observation.
$ws.Range("A24").Value2 = 9
$ws.Range("B24").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value2 = "Metropolitana"
$ws.Range("D24").Value2 = 45125
$ws.Range("E24").Value2 = 13
$ws.Range("F24").Value2 = 100112010
$ws.Range("G24").Value2 = "Achicoria"
$ws.Range("H24").Value2 = "Sin especificar"
$ws.Range("I24").Value2 = "Primera"
$ws.Range("J24").Value2 = 70
$ws.Range("K24").Value2 = 7000
$ws.Range("L24").Value2 = 7000
$ws.Range("M24").Value2 = 7000
$ws.Range("N24").Value2 = "`$/caja 16 unidades"
$ws.Range("O24").Value2 = "Provincia de Quillota"
$ws.Range("P24").Value2 = 438
$ws.Range("Q24").Value2 = 16
$ws.Range("R24").Value2 = "Hortaliza"
